$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (quantities) first so the shared-strings table ends up ordered
# the same way Excel produced it when this feature was added.
$ws.Range("C2").Value = "Name1"
$ws.Range("D2").Value = "Name2"
$ws.Range("E2").Value = "Name1"
$ws.Range("F2").Value = "Name2"
$ws.Range("G2").Value = "Name3"
$ws.Range("H2").Value = "Name1"
$ws.Range("I2").Value = "Name1"
$ws.Range("J2").Value = "Name1"
$ws.Range("K2").Value = "Name1"
$ws.Range("L2").Value = "Name1"

# Row 1 (quantity header list)
$ws.Range("C1").Value = "Name1"
$ws.Range("D1").Value = "Name3"
$ws.Range("E1").Value = "Name3"
$ws.Range("F1").Value = "Name1"
$ws.Range("G1").Value = "Name1"
$ws.Range("H1").Value = "Name1"
$ws.Range("I1").Value = "Name1"
$ws.Range("J1").Value = "Name1"
$ws.Range("K1").Value = "Name1"
$ws.Range("L1").Value = "Name1"

$ws.Range("L1").Select()
